$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.455.33'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.827.69'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.27'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5132'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3924'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07665'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.36%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.72'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.21%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.110'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("E12").Value = '  +1.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.293'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.25%  '
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.540'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.826.76'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.62'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +4.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001103'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06679'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.67'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("E21").Value = '  +0.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.150'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.492.21'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.254'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +8.38%  '
$ws.Range("E26").Value = '  +1.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.86'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.037.25'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.406'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.64'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.114'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1083'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.667'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("E34").Value = '  -0.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07021'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2210'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.920'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.19%  '
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("E39").Value = '  -1.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6263'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("E41").Value = '  -0.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.178'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.45%  '
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("E44").Value = '  -1.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.38'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5902'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.707'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.26'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.978'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.198'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.80%  '
$ws.Range("E51").Value = '  +0.55%  '
